$p = $ppt.ActivePresentation

# --- Change 1: Notes Master date placeholder "2016-07-01" -> "2016-07-02" ---
$nm = $p.NotesMaster
$dt = $nm.HeadersFooters.DateAndTime
$dt.Value = "2016-07-02"

# --- Change 2: Slide 1 subtitle placeholder gets an explicit position/size ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitle.Left = 164.38776
$subtitle.Top = 355.04886
$subtitle.Width = 391.22461
$subtitle.Height = 88.95122
